# Apply the "vendas1.xlsx" update:
#  - extend the "Vendas" used range from row 543 to row 554
#  - patch a handful of existing rows (mostly the "Data_recebida" / L column,
#    plus a few value corrections)
#  - append 11 new sales rows (544-554)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns actually populated on a "sales row" in this sheet (C, F, G are
# intentionally skipped throughout the data range).
$cols = @("A","B","D","E","H","I","J","K","L","M","N")

# ---------------------------------------------------------------------------
# 1. Row 426 - the "Data_recebida" (L) cell doesn't exist yet; clone the
#    number format from the nearest populated L cell (L425, same column,
#    same date style) before writing the value so it lands with style s="1"
#    like every other date cell instead of the default General style.
# ---------------------------------------------------------------------------
$ws.Range("L425").Copy()
$ws.Range("L426").PasteSpecial(-4122)
$ws.Range("L426").Value = 44681

# ---------------------------------------------------------------------------
# 2. Simple date-only corrections on existing rows (cells already exist with
#    the date style applied, so a plain value assignment is enough).
# ---------------------------------------------------------------------------
$ws.Range("L505").Value = 44681
$ws.Range("L509").Value = 44681
$ws.Range("L530").Value = 44681
$ws.Range("L536").Value = 44738
$ws.Range("L538").Value = 44681
$ws.Range("L539").Value = 44681
$ws.Range("L542").Value = 44708

# ---------------------------------------------------------------------------
# 3. Row 526 - value + received-value corrections
# ---------------------------------------------------------------------------
$ws.Range("H526").Value = 42420.51
$ws.Range("I526").Value = 42420.51
$ws.Range("L526").Value = 44926

# ---------------------------------------------------------------------------
# 4. Row 533 - now paid ("Pago" = S) with a received value
# ---------------------------------------------------------------------------
$ws.Range("J533").Value = "S"
$ws.Range("K533").Value = 157480
$ws.Range("L533").Value = 44677

# ---------------------------------------------------------------------------
# 5. Row 534 - received value + date
# ---------------------------------------------------------------------------
$ws.Range("K534").Value = 64450
$ws.Range("L534").Value = 44670

# ---------------------------------------------------------------------------
# 6. Row 535 - received value + date
# ---------------------------------------------------------------------------
$ws.Range("K535").Value = 30000
$ws.Range("L535").Value = 44670

# ---------------------------------------------------------------------------
# 7. Row 540 - description + values + paid flag + date
# ---------------------------------------------------------------------------
$ws.Range("E540").Value = "2 BALCÕES CASTAS E 4 CAIXOTES"
$ws.Range("H540").Value = 3498.81
$ws.Range("I540").Value = 3498.81
$ws.Range("J540").Value = "S"
$ws.Range("K540").Value = 3498.81
$ws.Range("L540").Value = 44650

# ---------------------------------------------------------------------------
# 8. Row 541 - values + date
# ---------------------------------------------------------------------------
$ws.Range("H541").Value = 4500
$ws.Range("I541").Value = 4500
$ws.Range("K541").Value = 4500
$ws.Range("L541").Value = 44655

# ---------------------------------------------------------------------------
# 9. New rows 544-554 - clone the row-543 formatting (per populated column,
#    so we don't manufacture cells in the unused C/F/G columns) into each new
#    row, then fill in the values.
# ---------------------------------------------------------------------------
$newRows = @(
    @{ Row=544; A=543; B=44644; D="V3A"; E="PETROBRAS ROVR NA RIO2C"; H=269000; I=269000; J="N"; K=0;        L=44718 },
    @{ Row=545; A=544; B=44655; D="ÓTIMA CONCESSIONÁRIA"; E="WARNER ANIMAIS FANTÁSTICOS"; H=83000; I=83000; J="N"; K=83000; L=44724 },
    @{ Row=546; A=545; B=44655; D="NETZA"; E="STAND EVE RIO2C"; H=64459.99; I=64459.99; J="N"; K=64459.99; L=44718 },
    @{ Row=547; A=546; B=44658; D="CERVEJARIA  PRAYA"; E="BAILE DO ENCANTO"; H=50000; I=50000; J="N"; K=50000; L=44681 },
    @{ Row=548; A=547; B=44658; D="AGÊNCIA TERRUÁ"; E="STAND BANCO DO BRASIL RIO2C"; H=200000; I=200000; J="N"; K=200000; L=44719 },
    @{ Row=549; A=548; B=44664; D="LVHM"; E="REMONTAGEM BAR BELVERDE"; H=7500; I=7500; J="N"; K=7500; L=44681 },
    @{ Row=550; A=549; B=44665; D="GLOBO SAT"; E="TÚNEL DO AMOR"; H=300000; I=300000; J="N"; K=300000; L=44712 },
    @{ Row=551; A=550; B=44666; D="GLOBO COMUNICAÇÕES"; E="CUBO RIO2C"; H=75000; I=75000; J="N"; K=75000; L=44712 },
    @{ Row=552; A=551; B=44667; D="BE COMUNICA"; E="STAND DASA"; H=85248; I=85248; J="N"; K=85248; L=44717 },
    @{ Row=553; A=552; B=44683; D="DIALOGO URBANO - BARBARA SOLEDADE"; E="STAND SUBMARINO (OBVIOUS)"; H=20000; I=20000; J="N"; K=20000; L=44742 },
    @{ Row=554; A=553; B=44685; D="VOID"; E="VOID TIJUCA"; H=49500; I=49500; J="N"; K=49500; L=44713 }
)

foreach ($r in $newRows) {
    $n = $r.Row
    $prev = $n - 1

    foreach ($col in $cols) {
        $ws.Range($col + $prev).Copy()
        $ws.Range($col + $n).PasteSpecial(-4122)
    }

    $ws.Range("A$n").Value = $r.A
    $ws.Range("B$n").Value = $r.B
    $ws.Range("D$n").Value = $r.D
    $ws.Range("E$n").Value = $r.E
    $ws.Range("H$n").Value = $r.H
    $ws.Range("I$n").Value = $r.I
    $ws.Range("J$n").Value = $r.J
    $ws.Range("K$n").Value = $r.K
    $ws.Range("L$n").Value = $r.L
    $ws.Range("M$n").Value = 0
    $ws.Range("N$n").Value = 0
}

# ---------------------------------------------------------------------------
# 10. Extend the "Vendas" defined name to the new used range
# ---------------------------------------------------------------------------
$wb.Names.Item("Vendas").RefersTo = "='Vendas'!`$A`$1:`$N`$554"
